# System documentation update:
#  1. "classes/methods" -> "classes/methods/interface" in the Javadoc bullet.
#  2. "a directory names "docs"." -> "a directory named "docs"." (fixing the
#     grammar in the demo javadoc comment; this also clears the stale
#     w:proofErr grammar-check markers that surrounded "directory names").

$d = $word.ActiveDocument

$quote = [char]8220
$endQuote = [char]8221

# --- Edit 1: classes/methods -> classes/methods/interface -----------------
$oldText1 = "Javadoc comments are provided for all classes/methods to document their functionality and usage."
$newText1 = "Javadoc comments are provided for all classes/methods/interface to document their functionality and usage."

$r1 = $d.Content
$r1.Find.Execute($oldText1, $true, $false, $false, $false, $false, $true, 1, $false, $newText1, 2)

# --- Edit 2: "a directory names "docs"." -> "a directory named "docs"." ---
$oldText2 = "This will generate documentation for all Java files in the current directory, and its subdirectories, and will output the HTML files to a directory names " + $quote + "docs" + $endQuote + "."
$newText2 = "This will generate documentation for all Java files in the current directory, and its subdirectories, and will output the HTML files to a directory named " + $quote + "docs" + $endQuote + "."

$r2 = $d.Content
$r2.Find.Execute($oldText2, $true, $false, $false, $false, $false, $true, 1, $false, $newText2, 2)
